# DutyList.xlsx update: add "dutyLoad" (numeric, 0.00) and "dutyType" (text)
# columns to the duty table, replacing the old "taskCount" column.
#
# Target layout (Sheet1, A1:H4):
#   Row1 (header, bold): Name | Region | startDate | startTime | endDate | endTime | dutyLoad | dutyType
#   Row2: deneme1 | deneme | 2020-02-15 | 10:00 | 2020-02-15 | 14:00 | 10.23 | S
#   Row3: deneme2 | deneme | 2020-02-15 | 10:00 | 2020-02-15 | 14:00 | 20.50 | A
#   Row4: deneme3 | deneme | 2020-02-16 | 02:00 | 2020-02-16 | 08:00 | 12.60 | SA
#
# Columns A-F are untouched by this edit; only G (repurposed from the old
# "taskCount" integer column into a "dutyLoad" decimal column) and the new
# H ("dutyType") column change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# --- Write all the new values first, column G before column H, top to
#     bottom. This keeps new shared-string entries in the same relative
#     order Excel produced them in (dutyLoad, dutyType, S, A, SA). ---
$ws.Range("G1").Value = "dutyLoad"
$ws.Range("H1").Value = "dutyType"
$ws.Range("G2").Value = 10.23
$ws.Range("H2").Value = "S"
$ws.Range("G3").Value = 20.5
$ws.Range("H3").Value = "A"
$ws.Range("G4").Value = 12.6
$ws.Range("H4").Value = "SA"

# --- Formatting. H1 before G1 so the two brand-new header cell styles are
#     created in the same relative order as the target workbook (H1's
#     "bold + horizontal-center only" style, then G1's "bold + 0.00 +
#     center/center" style). ---
$ws.Range("H1").HorizontalAlignment = $xlCenter
$ws.Range("H1").Font.Bold = $true

$ws.Range("G1").HorizontalAlignment = $xlCenter
$ws.Range("G1").VerticalAlignment = $xlCenter
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").NumberFormat = "0.00"

# dutyLoad values: centered, two-decimal number format.
$ws.Range("G2:G4").HorizontalAlignment = $xlCenter
$ws.Range("G2:G4").VerticalAlignment = $xlCenter
$ws.Range("G2:G4").NumberFormat = "0.00"

# dutyType values: centered text, same look as the other text columns (A,C-F).
$ws.Range("H2:H4").HorizontalAlignment = $xlCenter
$ws.Range("H2:H4").VerticalAlignment = $xlCenter
$ws.Range("H2:H4").NumberFormat = "@"
